$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns to snake_case field names
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'

# Title-case the connector words ("de", "del", "la", "las", "el", "los", "y")
# in state/municipality names for consistent casing
$ws.Range('B6').Value = 'Pabellón De Arteaga'
$ws.Range('B7').Value = 'Rincón De Romos'
$ws.Range('B8').Value = 'San José De Gracia'
$ws.Range('B12').Value = 'Playas De Rosarito'
$ws.Range('B25').Value = 'Amatenango De La Frontera'
$ws.Range('B28').Value = 'Bejucal De Ocampo'
$ws.Range('B34').Value = 'Chiapa De Corzo'
$ws.Range('B39').Value = 'Comitán De Domínguez'
$ws.Range('B48').Value = 'Marqués De Comillas'
$ws.Range('B49').Value = 'Mazapa De Madero'
$ws.Range('B54').Value = 'Ocozocoautla De Espinosa'
$ws.Range('B62').Value = 'Salto De Agua'
$ws.Range('B63').Value = 'San Cristóbal De Las Casas'
$ws.Range('B86').Value = 'Hidalgo Del Parral'
$ws.Range('A92').Value = 'Ciudad De México'
$ws.Range('B95').Value = 'Cuajimalpa De Morelos'
$ws.Range('A110').Value = 'Coahuila De Zaragoza'
$ws.Range('B123').Value = 'Villa De Álvarez'
$ws.Range('B131').Value = 'Nombre De Dios'
$ws.Range('B136').Value = 'San Juan Del Río'
$ws.Range('A143').Value = 'Estado De México'
$ws.Range('B143').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B145').Value = 'Almoloya De Alquisiras'
$ws.Range('B147').Value = 'Atizapán De Zaragoza'
$ws.Range('B156').Value = 'Ecatepec De Morelos'
$ws.Range('B162').Value = 'Ixtapan De La Sal'
$ws.Range('B163').Value = 'Ixtapan Del Oro'
$ws.Range('B171').Value = 'Naucalpan De Juárez'
$ws.Range('B175').Value = 'San Felipe Del Progreso'
$ws.Range('B176').Value = 'San José Del Rincón'
$ws.Range('B178').Value = 'San Simón De Guerrero'
$ws.Range('B182').Value = 'Tenango Del Valle'
$ws.Range('B185').Value = 'Tlalnepantla De Baz'
$ws.Range('B188').Value = 'Valle De Bravo'
$ws.Range('B191').Value = 'Villa De Allende'
$ws.Range('B204').Value = 'Jaral Del Progreso'
$ws.Range('B214').Value = 'San Diego De La Unión'
$ws.Range('B216').Value = 'San Francisco Del Rincón'
$ws.Range('B217').Value = 'San Miguel De Allende'
$ws.Range('B218').Value = 'Silao De La Victoria'
$ws.Range('B222').Value = 'Valle De Santiago'
$ws.Range('B228').Value = 'Acapulco De Juárez'
$ws.Range('B230').Value = 'Alcozauca De Guerrero'
$ws.Range('B234').Value = 'Atoyac De Álvarez'
$ws.Range('B235').Value = 'Ayutla De Los Libres'
$ws.Range('B238').Value = 'Buenavista De Cuéllar'
$ws.Range('B239').Value = 'Chilapa De Álvarez'
$ws.Range('B240').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B241').Value = 'Coahuayutla De José María Izazaga'
$ws.Range('B242').Value = 'Cochoapa El Grande'
$ws.Range('B247').Value = 'Coyuca De Benítez'
$ws.Range('B248').Value = 'Coyuca De Catalán'
$ws.Range('B251').Value = 'Cuetzala Del Progreso'
$ws.Range('B252').Value = 'Cutzamala De Pinzón'
$ws.Range('B258').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B259').Value = 'Iguala De La Independencia'
$ws.Range('B263').Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range('B277').Value = 'Taxco De Alarcón'
$ws.Range('B280').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B281').Value = 'Tixtla De Guerrero'
$ws.Range('B284').Value = 'Tlalixtaquilla De Maldonado'
$ws.Range('B285').Value = 'Tlapa De Comonfort'
$ws.Range('B287').Value = 'Técpan De Galeana'
$ws.Range('B290').Value = 'Zihuatanejo De Azueta'
$ws.Range('B299').Value = 'Atotonilco El Grande'
$ws.Range('B302').Value = 'Cuautepec De Hinojosa'
$ws.Range('B306').Value = 'Huejutla De Reyes'
$ws.Range('B311').Value = 'Mineral Del Monte'
$ws.Range('B312').Value = 'Mixquiahuala De Juárez'
$ws.Range('B314').Value = 'Pachuca De Soto'
$ws.Range('B318').Value = 'Santiago De Anaya'
$ws.Range('B321').Value = 'Tenango De Doria'
$ws.Range('B323').Value = 'Tepeji Del Río De Ocampo'
$ws.Range('B324').Value = 'Tula De Allende'
$ws.Range('B325').Value = 'Tulancingo De Bravo'
$ws.Range('B326').Value = 'Zacualtipán De Ángeles'
$ws.Range('B330').Value = 'Ahualulco De Mercado'
$ws.Range('B334').Value = 'Atemajac De Brizuela'
$ws.Range('B336').Value = 'Atotonilco El Alto'
$ws.Range('B337').Value = 'Autlán De Navarro'
$ws.Range('B346').Value = 'Cuautitlán De García Barragán'
$ws.Range('B355').Value = 'Encarnación De Díaz'
$ws.Range('B358').Value = 'Huejuquilla El Alto'
$ws.Range('B359').Value = 'Ixtlahuacán De Los Membrillos'
$ws.Range('B363').Value = 'Jilotlán De Los Dolores'
$ws.Range('B367').Value = 'La Manzanilla De La Paz'
$ws.Range('B368').Value = 'Lagos De Moreno'
$ws.Range('B377').Value = 'San Juan De Los Lagos'
$ws.Range('B380').Value = 'San Miguel El Alto'
$ws.Range('B382').Value = 'Santa María Del Oro'
$ws.Range('B384').Value = 'Talpa De Allende'
$ws.Range('B385').Value = 'Tamazula De Gordiano'
$ws.Range('B389').Value = 'Tepatitlán De Morelos'
$ws.Range('B392').Value = 'Tizapán El Alto'
$ws.Range('B397').Value = 'Unión De Tula'
$ws.Range('B398').Value = 'Valle De Guadalupe'
$ws.Range('B401').Value = 'Yahualica De González Gallo'
$ws.Range('B402').Value = 'Zacoalco De Torres'
$ws.Range('B405').Value = 'Zapotitlán De Vadillo'
$ws.Range('B406').Value = 'Zapotlán Del Rey'
$ws.Range('B407').Value = 'Zapotlán El Grande'
$ws.Range('A409').Value = 'Michoacán De Ocampo'
$ws.Range('B427').Value = 'Coalcomán De Vázquez Pallares'
$ws.Range('B428').Value = 'Cojumatlán De Régules'
$ws.Range('B505').Value = 'Coatlán Del Río'
$ws.Range('B511').Value = 'Puente De Ixtla'
$ws.Range('B515').Value = 'Tlaltizapán De Zapata'
$ws.Range('B525').Value = 'Bahía De Banderas'
$ws.Range('B531').Value = 'Santa María Del Oro'
$ws.Range('B542').Value = 'San Nicolás De Los Garza'
$ws.Range('B545').Value = 'Acatlán De Pérez Figueroa'
$ws.Range('B548').Value = 'Ayoquezco De Aldama'
$ws.Range('B551').Value = 'Chalcatongo De Hidalgo'
$ws.Range('B553').Value = 'Coicoyán De Las Flores'
$ws.Range('B554').Value = 'Constancia Del Rosario'
$ws.Range('B556').Value = 'Cuilápam De Guerrero'
$ws.Range('B558').Value = 'Eloxochitlán De Flores Magón'
$ws.Range('B559').Value = 'Guadalupe De Ramírez'
$ws.Range('B560').Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range('B561').Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range('B562').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B563').Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range('B564').Value = 'Huajuapan De León'
$ws.Range('B565').Value = 'Huautla De Jiménez'
$ws.Range('B567').Value = 'Ixtlán De Juárez'
$ws.Range('B570').Value = 'Mariscala De Juárez'
$ws.Range('B572').Value = 'Mazatlán Villa De Flores'
$ws.Range('B574').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B575').Value = 'Oaxaca De Juárez'
$ws.Range('B576').Value = 'Ocotlán De Morelos'
$ws.Range('B577').Value = 'Putla Villa De Guerrero'
$ws.Range('B582').Value = 'San Agustín De Las Juntas'
$ws.Range('B585').Value = 'San Antonino El Alto'
$ws.Range('B596').Value = 'San Felipe Jalapa De Díaz'
$ws.Range('B603').Value = 'San Francisco Del Mar'
$ws.Range('B626').Value = 'San Juan Del Río'
$ws.Range('B647').Value = 'San Miguel Del Puerto'
$ws.Range('B648').Value = 'San Miguel El Grande'
$ws.Range('B653').Value = 'San Pablo Villa De Mitla'
$ws.Range('B662').Value = 'San Pedro El Alto'
$ws.Range('B663').Value = 'San Pedro Y San Pablo Ayutla'
$ws.Range('B677').Value = 'Santa Cruz Tacache De Mina'
$ws.Range('B680').Value = 'Santa Cruz De Bravo'
$ws.Range('B681').Value = 'Santa Inés Del Monte'
$ws.Range('B692').Value = 'Santa María Jalapa Del Marqués'
$ws.Range('B732').Value = 'Tamazulápam Del Espíritu Santo'
$ws.Range('B733').Value = 'Tataltepec De Valdés'
$ws.Range('B734').Value = 'Tezoatlán De Segura Y Luna'
$ws.Range('B735').Value = 'Tlacolula De Matamoros'
$ws.Range('B736').Value = 'Totontepec Villa De Morelos'
$ws.Range('B738').Value = 'Villa Sola De Vega'
$ws.Range('B739').Value = 'Villa Talea De Castro'
$ws.Range('B740').Value = 'Villa De Etla'
$ws.Range('B741').Value = 'Villa De Tamazulápam Del Progreso'
$ws.Range('B742').Value = 'Villa De Tututepec'
$ws.Range('B743').Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range('B744').Value = 'Villa De Zaachila'
$ws.Range('B746').Value = 'Zimatlán De Álvarez'
$ws.Range('B757').Value = 'Chalchicomula De Sesma'
$ws.Range('B781').Value = 'Huehuetlán El Chico'
$ws.Range('B782').Value = 'Huehuetlán El Grande'
$ws.Range('B787').Value = 'Ixcamilpa De Guerrero'
$ws.Range('B789').Value = 'Izúcar De Matamoros'
$ws.Range('B797').Value = 'Los Reyes De Juárez'
$ws.Range('B803').Value = 'Palmar De Bravo'
$ws.Range('B813').Value = 'San Nicolás De Los Ranchos'
$ws.Range('B815').Value = 'San Salvador El Verde'
$ws.Range('B821').Value = 'Tecali De Herrera'
$ws.Range('B825').Value = 'Tepanco De López'
$ws.Range('B826').Value = 'Tepango De Rodríguez'
$ws.Range('B829').Value = 'Tepexi De Rodríguez'
$ws.Range('B830').Value = 'Tetela De Ocampo'
$ws.Range('B833').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B853').Value = 'Amealco De Bonfil'
$ws.Range('B854').Value = 'Cadereyta De Montes'
$ws.Range('B857').Value = 'Landa De Matamoros'
$ws.Range('B860').Value = 'Pinal De Amoles'
$ws.Range('B862').Value = 'San Juan Del Río'
$ws.Range('B878').Value = 'Soledad De Graciano Sánchez'
$ws.Range('B882').Value = 'Villa De Ramos'
$ws.Range('B917').Value = 'San Miguel De Horcasitas'
$ws.Range('B942').Value = 'Contla De Juan Cuamatzi'
$ws.Range('B945').Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range('B948').Value = 'San Pablo Del Monte'
$ws.Range('B951').Value = 'Tetla De La Solidaridad'
$ws.Range('A955').Value = 'Veracruz De Ignacio De La Llave'
$ws.Range('B960').Value = 'Amatlán De Los Reyes'
$ws.Range('B974').Value = 'Cosamaloapan De Carpio'
$ws.Range('B989').Value = 'Juchique De Ferrer'
$ws.Range('B993').Value = 'Martínez De La Torre'
$ws.Range('B997').Value = 'Mixtla De Altamirano'
$ws.Range('B999').Value = 'Nanchital De Lázaro Cárdenas Del Río'
$ws.Range('B1004').Value = 'Paso De Ovejas'
$ws.Range('B1006').Value = 'Poza Rica De Hidalgo'
$ws.Range('B1013').Value = 'Sayula De Alemán'
$ws.Range('B1015').Value = 'Soledad De Doblado'
$ws.Range('B1017').Value = 'Tatahuicapan De Juárez'
$ws.Range('B1029').Value = 'Vega De Alatorre'
$ws.Range('B1055').Value = 'Concepción Del Oro'
$ws.Range('B1065').Value = 'Mezquital Del Oro'
$ws.Range('B1066').Value = 'Nochistlán De Mejía'
$ws.Range('B1071').Value = 'Teúl De González Ortega'
$ws.Range('B1072').Value = 'Tlaltenango De Sánchez Román'
$ws.Range('B1073').Value = 'Trinidad García De La Cadena'
$ws.Range('B1078').Value = 'Villa De Cos'

# Remove trailing metadata/footer rows (sample size, source, author, date)
$ws.Rows("1084:1088").Delete()
